$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.024.62'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +3.08%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.599.86'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +2.96%  '

$ws.Range("E4").Value = '  -0.08%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '212.24'
$c.Style = "Normal"

$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("E7").Value = '  +1.56%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.247'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +2.06%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.0613'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.83%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '18.06'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.96%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0815'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +4.31%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.825.56'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.16%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.602.04'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +3.12%  '

$ws.Range("E14").Value = '  +0.31%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.509'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.92%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '26.025.69'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +3.13%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '60.26'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.41%  '

$ws.Range("E18").Value = '  +1.94%  '

$ws.Range("E19").Value = '  -0.09%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '201.11'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +8.44%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '4.21'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.64%  '

$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("E23").Value = '  +2.80%  '

$ws.Range("E24").Value = '  +12.88%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '141.29'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '

$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("E27").Value = '  -5.79%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '15.15'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.12%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '6.41'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("E30").Value = '  +1.69%  '

$ws.Range("E31").Value = '  +1.54%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.10'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +2.63%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.95'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.30%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.47'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.78%  '

$ws.Range("E35").Value = '  +1.09%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0165'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +11.08%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.124.23'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +3.92%  '

$ws.Range("E38").Value = '  -0.05%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.787'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.85%  '

$ws.Range("E40").Value = '  +2.39%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.490'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.63%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.782'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.56%  '

$ws.Range("B43").Value = 'RocketPoolETH'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.737.49'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.15%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.14'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.96%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '93.10'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.68%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.50'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.87%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '53.28'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.83%  '

$ws.Range("E48").Value = '  -0.08%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.409'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.09%  '

$ws.Range("E50").Value = '  +0.14%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0₇0924'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -16.24%  '
